$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6091.2856
$ws.Range("J70").Value = 5606.5
$ws.Range("L70").Value = 16819.5
$ws.Range("N70").Value = -17359.5
$ws.Range("H73").Value = 6091.2856
$ws.Range("J73").Value = 5606.5
$ws.Range("L73").Value = 16819.5
$ws.Range("N73").Value = -18691.5
$ws.Range("H76").Value = 3898.8333
$ws.Range("I76").Value = 3650
$ws.Range("K76").Value = 3650
$ws.Range("M76").Value = -3335
$ws.Range("H79").Value = 3898.8333
$ws.Range("I79").Value = 3650
$ws.Range("K79").Value = 3650
$ws.Range("M79").Value = -2558
$ws.Range("H113").Value = 2499.25
$ws.Range("I113").Value = 998.5
$ws.Range("K113").Value = 998.5
$ws.Range("M113").Value = 2255.5
$ws.Range("H116").Value = 25761696
$ws.Range("I116").Value = 26988070
$ws.Range("K116").Value = 26988070
$ws.Range("M116").Value = -26984628
$ws.Range("H129").Value = 4847.5
$ws.Range("I129").Value = 1097
$ws.Range("K129").Value = 3291
$ws.Range("M129").Value = 1709
$ws.Range("H137").Value = 11447699
$ws.Range("I137").Value = 628166.7
$ws.Range("K137").Value = 1884500.1
$ws.Range("M137").Value = -1881950.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 950
$ws.Range("I2").Value = 950
$ws.Range("K2").Value = 950
$ws.Range("M2").Value = -837
$ws.Range("H61").Value = 4778.1055
$ws.Range("I61").Value = 2837.2307
$ws.Range("K61").Value = 2837.2307
$ws.Range("M61").Value = -2625.2307
$ws.Range("I74").Value = 20834144
$ws.Range("K74").Value = 20834144
$ws.Range("M74").Value = -20833270
$ws.Range("I77").Value = 20834144
$ws.Range("K77").Value = 104170720
$ws.Range("M77").Value = -104166352
$ws.Range("H102").Value = 764.2857
$ws.Range("I102").Value = 764.2857
$ws.Range("K102").Value = 764.2857
$ws.Range("M102").Value = 857.7143
$ws.Range("H110").Value = 859.75
$ws.Range("I110").Value = 528.3333
$ws.Range("K110").Value = 528.3333
$ws.Range("M110").Value = 1516.6667
$ws.Range("H116").Value = 950
$ws.Range("I116").Value = 950
$ws.Range("K116").Value = 950
$ws.Range("M116").Value = 1344
$ws.Range("H122").Value = 5303.3125
$ws.Range("I122").Value = 3857.375
$ws.Range("K122").Value = 11572.125
$ws.Range("M122").Value = -9122.125
$ws.Range("H132").Value = 18751.371
$ws.Range("I132").Value = 20296.709
$ws.Range("K132").Value = 60890.12699999999
$ws.Range("M132").Value = -58360.12699999999
$ws.Range("H136").Value = 4778.1055
$ws.Range("I136").Value = 2837.2307
$ws.Range("K136").Value = 8511.6921
$ws.Range("M136").Value = -5961.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 950
$ws.Range("I3").Value = 950
$ws.Range("K3").Value = 950
$ws.Range("M3").Value = -836
$ws.Range("H132").Value = 79537.5
$ws.Range("J132").Value = 79537.5
$ws.Range("L132").Value = 79537.5
$ws.Range("N132").Value = -89657.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17547110
$ws.Range("I31").Value = 19233586
$ws.Range("K31").Value = 19233586
$ws.Range("M31").Value = -19233291
$ws.Range("H34").Value = 17547110
$ws.Range("I34").Value = 19233586
$ws.Range("K34").Value = 19233586
$ws.Range("M34").Value = -19233384
$ws.Range("H58").Value = 3668.875
$ws.Range("I58").Value = 3566
$ws.Range("K58").Value = 3566
$ws.Range("M58").Value = -3363
$ws.Range("H99").Value = 4251.1055
$ws.Range("I99").Value = 2077.4
$ws.Range("J99").Value = 6666.3335
$ws.Range("K99").Value = 2077.4
$ws.Range("L99").Value = 6666.3335
$ws.Range("M99").Value = -579.4000000000001
$ws.Range("N99").Value = -9662.333500000001
$ws.Range("H107").Value = 819.2069
$ws.Range("J107").Value = 914.8095
$ws.Range("L107").Value = 914.8095
$ws.Range("N107").Value = -4754.8095
$ws.Range("H126").Value = 4251.1055
$ws.Range("I126").Value = 2077.4
$ws.Range("J126").Value = 6666.3335
$ws.Range("K126").Value = 6232.200000000001
$ws.Range("L126").Value = 19999.0005
$ws.Range("M126").Value = -3762.200000000001
$ws.Range("N126").Value = -24939.0005
$ws.Range("H132").Value = 35089324
$ws.Range("I132").Value = 45978430
$ws.Range("J132").Value = 2198
$ws.Range("K132").Value = 137935290
$ws.Range("L132").Value = 6594
$ws.Range("M132").Value = -137932760
$ws.Range("N132").Value = -11654
$ws.Range("H134").Value = 3227.111
$ws.Range("I134").Value = 2546.2307
$ws.Range("K134").Value = 7638.6921
$ws.Range("M134").Value = -5103.6921
$ws.Range("H136").Value = 3668.875
$ws.Range("I136").Value = 3566
$ws.Range("K136").Value = 10698
$ws.Range("M136").Value = -8148
$ws.Range("H141").Value = 163134.2
$ws.Range("J141").Value = 191396.5
$ws.Range("L141").Value = 191396.5
$ws.Range("N141").Value = -201756.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1048.4615
$ws.Range("J121").Value = 1273.9
$ws.Range("L121").Value = 3821.7
$ws.Range("N121").Value = -6441.700000000001
$ws.Range("H134").Value = 7963.8887
$ws.Range("I134").Value = 5209.375
$ws.Range("K134").Value = 15628.125
$ws.Range("M134").Value = -10558.125
$ws.Range("H138").Value = 3181.1667
$ws.Range("I138").Value = 3264.2666
$ws.Range("J138").Value = 2765.6667
$ws.Range("K138").Value = 9792.799800000001
$ws.Range("L138").Value = 8297.000100000001
$ws.Range("M138").Value = -4652.799800000001
$ws.Range("N138").Value = -18577.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58318.285
$ws.Range("I80").Value = 66233.69
$ws.Range("K80").Value = 66233.69
$ws.Range("M80").Value = -65235.69
$ws.Range("H83").Value = 58318.285
$ws.Range("I83").Value = 66233.69
$ws.Range("K83").Value = 331168.45
$ws.Range("M83").Value = -326176.45
$ws.Range("H122").Value = 307522.97
$ws.Range("I122").Value = 558283.4
$ws.Range("J122").Value = 6610.467
$ws.Range("K122").Value = 1674850.2
$ws.Range("L122").Value = 19831.401
$ws.Range("M122").Value = -1672400.2
$ws.Range("N122").Value = -24731.401
$ws.Range("H132").Value = 81796.39999999999
$ws.Range("I132").Value = 85150.414
$ws.Range("K132").Value = 255451.242
$ws.Range("M132").Value = -252921.242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6065.926
$ws.Range("I7").Value = 5514.2144
$ws.Range("K7").Value = 5514.2144
$ws.Range("M7").Value = -5402.2144
$ws.Range("H46").Value = 3881.4
$ws.Range("I46").Value = 1738.5555
$ws.Range("K46").Value = 1738.5555
$ws.Range("M46").Value = -1550.5555
$ws.Range("H82").Value = 4568.4287
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 5163.1665
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 5163.1665
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -5885.1665
$ws.Range("H85").Value = 4568.4287
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 5163.1665
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 5163.1665
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -7659.1665
$ws.Range("H93").Value = 2042.2778
$ws.Range("I93").Value = 1776.3
$ws.Range("J93").Value = 2374.75
$ws.Range("K93").Value = 1776.3
$ws.Range("L93").Value = 2374.75
$ws.Range("M93").Value = -528.3
$ws.Range("N93").Value = -4870.75
$ws.Range("H126").Value = 6065.926
$ws.Range("I126").Value = 5514.2144
$ws.Range("K126").Value = 16542.6432
$ws.Range("M126").Value = -14072.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4177.6333
$ws.Range("I100").Value = 575.6842
$ws.Range("J100").Value = 10399.182
$ws.Range("K100").Value = 1151.3684
$ws.Range("L100").Value = 20798.364
$ws.Range("M100").Value = -610.3684000000001
$ws.Range("N100").Value = -21880.364
$ws.Range("H103").Value = 33999.75
$ws.Range("J103").Value = 33999.75
$ws.Range("L103").Value = 33999.75
$ws.Range("N103").Value = -36343.75
$ws.Range("H126").Value = 37039436
$ws.Range("I126").Value = 43480492
$ws.Range("J126").Value = 3372.75
$ws.Range("K126").Value = 130441476
$ws.Range("L126").Value = 10118.25
$ws.Range("M126").Value = -130439006
$ws.Range("N126").Value = -15058.25
$ws.Range("H132").Value = 5217.9214
$ws.Range("I132").Value = 1138.9
$ws.Range("K132").Value = 3416.7
$ws.Range("M132").Value = -886.7000000000003
